$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Code Review 4 scores that were filled in for three of the four team members
$ws.Range("F4").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("F7").Value = 100

# Josh McQueen (row 6) still hasn't submitted Code Review 4 - highlight his
# row (F6:G6) in gray to flag it as outstanding
$ws.Range("F6:G6").Interior.Color = 10921638

# Leave the cursor where the user was last working
$ws.Range("H7").Select()
